$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (G1) onto the new H1 header
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
